# ManageNews page and test
# Adds a new "ManageNewsPage" worksheet at the end of the workbook,
# populated with two test-news strings, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the current last sheet, so it lands at the
# end of the tab order (sheetId 3 / rId3), and rename it.
$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws.Name = "ManageNewsPage"

# Populate the test data (adds two new shared strings).
$ws.Range("A1").Value = "This is a test news"
$ws.Range("A2").Value = "This is a test news 1"

# Make the new sheet the active/selected tab.
$ws.Activate()

# Mirror the author's final UI state: the whole sheet selected (as if
# "Select All" was pressed) on the new sheet.
$ws.Cells.Select()
